# Adds a new concluding paragraph to the "Conclusões do avaliador sobre o
# teste com o usuário" section, right after the "Sugestões de melhoria..."
# paragraph and before the trailing blank paragraph at the end of the body.

$d = $word.ActiveDocument

# Anchor on the unique tail of the "Sugestões de melhoria..." paragraph so
# this doesn't depend on fragile paragraph indices.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "reflitam exatamente o conteúdo que será exibido para o usuário.",
    $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor text not found"
}

# Collapse to the end of the match (= end of that paragraph, just before its
# paragraph mark) and insert a brand-new paragraph right after it.
$anchor.Collapse(0)  # wdCollapseEnd
$anchor.InsertParagraphAfter() | Out-Null

# Re-acquire the freshly inserted paragraph (now the second-to-last one,
# since the original trailing empty paragraph is still last). It already
# inherited the body-text formatting (Arial 10pt, double-spaced, pt-BR)
# from the paragraph mark it was split off from; just fill in the text.
$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count - 1)

$newPara.Range.Text = "Funcionalidades adicionais que podem otimizar a experiência do usuário são implementar caixas de diálogos/alertas com mensagens objetivas e claras para auxiliar a navegação do site, permitir ao usuário cancelar ou confirmar ações, fazer com que a logo do site redirecione o usuário para a página inicial."
